$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Phone Number"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Number of guests"

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Aamer"
$ws.Cells.Item(2, 3).Value = 966535288851
$ws.Cells.Item(2, 5).Value = 4

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Aya"
$ws.Cells.Item(3, 3).Value = 905525194276
$ws.Cells.Item(3, 5).Value = 2

# Apply General number format to column C
$ws.Range("C1:C3").NumberFormat = "General"

# Hyperlinks for email addresses (explicit display text keeps the cell's visible value)
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:aamermurhaf@hotmail.com", "", "", "aamermurhaf@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:ayashams42@yahoo.com", "", "", "ayashams42@yahoo.com")

# Column widths (closest achievable widths to the authored 29.453125 / 24.90625 / 15.6328125)
$ws.Range("C1").EntireColumn.ColumnWidth = 28.666666666666668
$ws.Range("D1").EntireColumn.ColumnWidth = 24
$ws.Range("E1").EntireColumn.ColumnWidth = 14.833333333333332

$ws.Range("G6").Select()
